$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for rows 108-114 (columns C=description, D=causes, E=symptoms, F=treatments, G=contraindications)
$data = @(
    @(108, "Anorexia nervosa is a condition marked by sever weight loss resulting from an unhealthy restricting of caloric intake. In conjunction with eating an inadequate amount of food, a person with anorexia may also try losing weight by taking diuretics, laxatives, and vomiting after eating . This can cause a severe lack of nutrition, resulting in the body thinning far beyond a healthy level.",
              "Anorexia si a pyschological disorder in which a person's perception of their weight is distorted. Contributing factors towards the development of anorexia include environmental influences that put an over-emphasis on being thin and pyschological issues such as obsessive compulsive disorder that make it easier to not eat by sticking to set goals.",
              "Symptoms are wide-ranging, from severe thinness and weight loss, to problems with tooth decay from excessive vomiting. A lack of proper nutrition can lead to fatigue, weakness, thinning hair, and dehydration. A person with anorexia may develop anemia, osteoporosis, abnormalities in hormone production and regulation, issues with kidney function, and muscle atrophy. Pyschologically, a person may withdraw socially, and try to hide their anorexia by wearing clothing to hide their weight loss. A person may become irritable and skip meals. When confronted, a person may deny skipping meals or lie about how much food they have eaten.",
              "Treatment for anorexia may include hospitalization in severe cases, where the body is not receiving adequate nutrition over a long period of time. While in the hospital, a person will be given fluids to balance hydration and electrolyte levels, and be treated for issues possibly relating to the heart, liver, and kidneys. A feeding tube may be inserted to ensure a person is receiving enough nutrients. Aside from a stay in the hospital, a person may be seen by a mental health professional to deal with underlying causes of anorexia. A person may work closely with a dietician to maintain a healthy diet.",
              "Massage therapy is not contraindicated for clients with anorexia. A lighter massage should be performed to prevent damage to bones and other structures in the body."),
    @(109, "Anxiety disorders are a group of disorders in where a person experiences an increased amount of anxiety, often at a level where the person feels they are in a life-or-death situation. These anxious moments can adversely affect a person's daily life and relationships. Different types of anxiety disorders include panic disorder, social anxiety disorder, and general anxiety disorder. A person with an anxiety disorder may also suffer from addiction, depression, chronic pain, or even attempt suicide.",
              "Often times, underlying medical issues may be the cause of anxiety disorders. Some examples include cardiovascular diseases, substance withdrawal, asthma, and hyperthyroidism. Medications may also induce anxiety.",
              "Symptoms very depending on the type of anxiety a person suffers from. General symptoms include nervousness, increased breathing and heart rate, shaking, trembling, digestive issues, and the inability to stop thinking about whatever it is that is putting a person into a state of anxiety. Panic anxiety may leave a person with a feeling of fear or terror, and lead to chest pain, shortness of breath, and a rapid heart beat. Social anxiety may lead a person to avoid contact with others as a means to protect themselves from situations they might find themselves being embarrassed or self-conscious in.",
              "Psychotherapy is the primary treatment for anxiety disorders. Cognitive Behavioral Therapy is used to help treat social anxiety by exposing a person to situations they may experience in real life, helping them learn how to cope with the situation when it arises. Medications may be prescribed, such antidepressants, that can help a person from feeling anxious.",
              "Anxiety disorders are not contraindicated for massage therapy."),
    @(110, "Bulimia is a psychological disorder in which a person ingests abnormally large portions of food, followed by purging of the food ingested. Purging of food can be caused by self-induced vomiting, by consuming laxatives or diuretics, and/or excessive exercise. This is all done in an effort to avoid weight gain. The use of self-induced vomiting and laxatives/diuretics is known as pruging bulimia. These involve ridding the body of the ingested food in some way that doesn't allow it to be properly digested. The use of excessive exercise is known as non-purging bulimia.",
              "The exact cause of bulimia is unknown, although there are several contributing factors that may lead to the development of the condition. Stress, histoy of abuse, trauma, low self-esteem, and having a negative image of one's body can all be contributing factors.",
              "Common symptoms may include dehydration, imbalances in electrolyte levels, fluctuating weight, lesions in the mouth due to excessive vomiting, chronic heart burun, and infertility. People with bulimia may exhibit some of the following traits and behaviours: frequent bathroom usage after eating, eating privately, smelling of vomit, and lacking control when eating.",
              "Bulimia is most commonly caused by low self-esteem and negative body image. Therefore, the primary treatment is therapy to help the patient overcome these pyschological issues.",
              "Bulimia is not contraindicated for massage therapy. Massage may help psychologically and improve self esteem and body image. A lighter massage may be indicated, as the client may have atrophied muscles and a weakended body due to lack of nutrients."),
    @(111, " Dementia is a general term for symptoms that involve loss of memory, social skills, and critical thinking skills. These skills may be impairred to the point where it affects a person's ability to function effectively on a daily basis. Dementia is caused by damage to neurons in the brain. There are numerous reasons for neurons to be damaged that can contribute to dementia.",
              "Dementia may be permanent, or may be reversed , depending on the cause. Causes of dementia that involve permanent damage to the brain include progressive disorders such as Alzheimer's disease. Strokes may permanently damage the brain. Similar to strokes, injuries to the brain such as consussions may lead to brain damage. Cases where dementia may be reversed include certain infections that cause high fever, autoimmune disorders such as multiple sclerosis, tumors in the brain, hypoxia, inadequate fluid intake, and hormonal problems such as having too much calcium in the body. As these conditions are treated, the dementia should go away.",
              " Dementia is most commonly associated with memory loss. This can lad to problems communicating, focusing, usign the body properly, and problem solving. A person may become depressed, develop anxiety or paranoia, and become easily irritated. ",
              "Dementia itself can't be cured. Treating the underlying cause is the best way to treat dementia. This can include medications to help brain functions, and physical therapy to help with any motor skills that may be affected.",
              "Dementia is not contraindicated for massage therapy."),
    @(112, "Depression is a disorder affecting mood, causing a person to feel sad, usually over a prolonged period of time. A person may become disinterested in normal daily activities, and become more distant. Depression usually develops during teenage years, but may occur at any time. If a person is experiencing depression, they should seek help from a mental health professional",
              "There are many different causes for depression, and no one person is the same as another. Causes include a change in the body's hormone leveles, such as during and after pregnancy, chagnes in the funciton of neurotransmitters in the brain, and even may be genetically passed down.",
              "Symptoms are wide-ranging, and people experience depression differently. Symptoms may include a feeling of sadness, fatigue, disinterest in eating, anxiety, insomnia, irritability, and constantly thinking about death and suicide. A person with depression may attempt suicide.",
              "If a person is diagnosed with depression, antidepressants known as selective serotonin reuptake inhibitors, or SSRIs are usually prescribed to help manage hormone levels and increase serotonin levels in the blood. Talking with a mental halth professional may help a person feel better. Exercise and a healthy diet may help with hormone imbalances and increase a person's self esteem, which can aid in lowering depression. Other relaxation techniques, such as massage therapy and yoga may be sought.",
              "Depression is not contraindicated for massage therapy. Massge may help relax a person with depression, increase self-esteem and body image, and help manage hormone leveles in the body."),
    @(113, "Insomnia is a disorder causing a lack of sufficient sleep. This may be the result of a person having difficulty fallling asleep, staying asleep, or waking too early. This can lead to many health complications.",
              "There are many reasons a person's sleep may be affected. Stress is the most common cause of insomnia, making it difficult to 'turn off' thoughts when its time to sleep. Having an irregular sleep schedule may make it difficult for the body to adjust when its time to sleep. Jet lag may also cause insomnia. and certian medications may interfere with normal sleep. Pain associated with medical conditions may also keep a person from sleeping.",
              "Common symptoms include trouble falling asleep, trouble staying asleep waking too early, general fatigue, irritablility, increased stress and anxiety, and depression. A person with insomnia may develop other health conditions as a result of not getting enough sleep, such as hypertension.",
              "Determining the root cause of the insomnia and treating that is the primary treatment for insomnia. Reducing stress, getting a person on a set sleep scheudle, increasing exercise, meditaiton, yoga, and massage therapy can all aid in eliminating insomnia. If these methods do not accomplish the goal, medications to help a person sleep may be prescribed. Over-the-counter sleep aids may help a person sleep, but should not be used long-term as a person may develop a dependency for them.",
              "Insomnia is not a contraindication for massage. Massage therapy may aid a person in reducing stress and helping them to obtain a normal sleep cycle."),
    @(114, "Veritgo is the senstion that a person's surroundings are moving or spinning, which may cause dizziness. This can be especially prevalent if a person is looking down from a tall height.",
              "Commonly, vertigo is the result of a problem with the inner ear. Menier's desease, which causes fluid buildup that changes pressure inside the ear, is a common cause. A buildup of calcium deposits in the inner ear can alter balance. Viral infections, such as vestibular neuritis, may also cause vertigo. Less commonly, tumors in the ear may cause vertigo.",
              "Verigo presents with a feeling of the environment around a person spinning, moving, swaying, or tilting. A person may become unbalanced, become nauseated, vomit, or develop headaches.",
              "Treatment for vertigo is usually dependent on the cause. If calcium is present in the inner ear, certain head and neck movements may be performed to aid the calcium in leaving the inner ear, allowing it to be broken down by the body. Medications may be prescribed to aid with nausea and fluid build up associated with Menier's disease. Surgery may be performed if there is a tumor present. Often times, no treatmetn is necessary, as the brain becomes acclimated to vertigo and the symptoms lessen or disappear.",
              "Vertigo is not contraindicated for massage. The massage therapist may need to help the client on and off the table to avoid accidental falls due to dizziness.")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

$ws.Range("E116").Select()
